$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("Username") and Column D ("Random") are grouped/merged so that
# rows sharing the same block take on a common value (grouping + merging
# per the commit message).
#
# Column B groups: rows 2-8 -> 6090, rows 9-15 -> 3182, rows 16-21 -> 1666
$ws.Range("B2:B8").Value = 6090
$ws.Range("B9:B15").Value = 3182
$ws.Range("B16:B21").Value = 1666

# Column D groups: rows 2-7 -> 3128, rows 8-15 -> 6962, rows 16-21 -> 1892
$ws.Range("D2:D7").Value = 3128
$ws.Range("D8:D15").Value = 6962
$ws.Range("D16:D21").Value = 1892

# Selected cell after the edit
$ws.Range("C12").Select()
